$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 441, shifting rows 441:527 down to 442:528
$ws.Rows.Item(441).Insert()

# Populate the newly inserted row 441 with the new record
$ws.Cells.Item(441, 1).Value = 3
$ws.Cells.Item(441, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(441, 3).Value = "Coquimbo"
$ws.Cells.Item(441, 4).Value = 44637
$ws.Cells.Item(441, 5).Value = 5
$ws.Cells.Item(441, 6).Value = 100114001
$ws.Cells.Item(441, 7).Value = "Papa"
$ws.Cells.Item(441, 8).Value = "Rosara"
$ws.Cells.Item(441, 9).Value = "1a (cosecha)"
$ws.Cells.Item(441, 10).Value = 2370
$ws.Cells.Item(441, 11).Value = 7500
$ws.Cells.Item(441, 12).Value = 8000
$ws.Cells.Item(441, 13).Value = 7946
$ws.Cells.Item(441, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(441, 15).Value = "Provincia de Talca"
$ws.Cells.Item(441, 16).Value = 318
$ws.Cells.Item(441, 17).Value = 25
$ws.Cells.Item(441, 18).Value = "Hortaliza"
